$d = $word.ActiveDocument

# --- Change 1: append a red, parenthetical note to the first paragraph ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1          # exclude the paragraph mark
$r1.Text = "This is a Microsoft word document.  "

$insPoint = $d.Range($r1.End, $r1.End)
$insPoint.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$insPoint.Font.Color = 255      # wdColorRed -> OOXML FF0000

$insPoint2 = $d.Range($insPoint.End, $insPoint.End)
$insPoint2.InsertAfter("rsion for main branch")
$insPoint2.Font.Color = 255

$insPoint3 = $d.Range($insPoint2.End, $insPoint2.End)
$insPoint3.InsertAfter(")")
$insPoint3.Font.Color = 255

# --- Change 2: remove the trailing "ank God almighty, we are free at last." paragraph ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.Delete()

# --- Change 3: drop the unused leftover styles (pasted-from-web cruft) ---
$namesToDelete = @("Heading 2", "Heading 4", "apple-converted-space", "Hyperlink",
                    "Heading 2 Char", "Heading 4 Char", "audio-tool", "subscribe",
                    "subscribe-more-info", "generic-title", "podcast-tools__subscribe-links")

# Resolve every target to its current (1-based) index first -- plain integers,
# not live object handles, so nothing goes stale once deletions start shifting
# the collection around.
$indexesToDelete = @()
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    if ($namesToDelete -contains $d.Styles.Item($i).NameLocal) {
        $indexesToDelete += $i
    }
}

# Delete highest index first so the indexes still to be visited never shift.
$indexesToDelete = $indexesToDelete | Sort-Object -Descending
foreach ($idx in $indexesToDelete) {
    $d.Styles.Item($idx).Delete()
}
